$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("SpPa", $false, $false, $false, $false, $false, $true, 1, $false, "PaSp", 2)
